$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.865.84'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.555.88'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.75'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.31'
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.555.01'
$ws.Range("E7").Value = '  -0.98%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("E9").Value = '  +5.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.133'
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.83'
$ws.Range("E11").Value = '  -2.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.413'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.154.55'
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000198'
$ws.Range("E14").Value = '  -5.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.18'
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.557.70'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.117'
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.685.44'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.04'
$ws.Range("E19").Value = '  -3.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.27'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.78'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '427.43'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.601'
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.90'
$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.694.53'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -4.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.07'
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.08'
$ws.Range("E30").Value = '  -3.03%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.560.10'
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.54'
$ws.Range("E34").Value = '  -3.46%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -7.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.67'
$ws.Range("E37").Value = '  -2.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.65'
$ws.Range("E38").Value = '  -3.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '177.31'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.32'
$ws.Range("E40").Value = '  -5.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0832'
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.05'
$ws.Range("E42").Value = '  -3.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.865'
$ws.Range("E43").Value = '  -3.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.63'
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.79'
$ws.Range("E45").Value = '  -5.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.15'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.45'
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("E50").Value = '  -4.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.922'
$ws.Range("E51").Value = '  -3.02%  '
